# Auto-generated edits applying cached value updates from the source diff.
# Each row's H-N columns (currentAveragePrice.. LeveProfitHQ) are refreshed
# with updated market-derived figures; some cells are cleared entirely where
# the source no longer emits a cached value for that column.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 817.2727
$ws.Range("I18").Value = 817.2727
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 817.2727
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -533.2727
$ws.Range("N18").ClearContents()

$ws.Range("H88").Value = 5641.75
$ws.Range("I88").Value = 511.77777
$ws.Range("J88").Value = 9839
$ws.Range("K88").Value = 511.77777
$ws.Range("L88").Value = 9839
$ws.Range("M88").Value = -105.77777
$ws.Range("N88").Value = -10651

$ws.Range("H91").Value = 5641.75
$ws.Range("I91").Value = 511.77777
$ws.Range("J91").Value = 9839
$ws.Range("K91").Value = 511.77777
$ws.Range("L91").Value = 9839
$ws.Range("M91").Value = 892.2222300000001
$ws.Range("N91").Value = -12647

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 886.1667
$ws.Range("I4").Value = 1060.4
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 1060.4
$ws.Range("L4").Value = 15
$ws.Range("M4").Value = -944.4000000000001
$ws.Range("N4").Value = -247

$ws.Range("H132").Value = 9402.32
$ws.Range("I132").Value = 5626.8335
$ws.Range("J132").Value = 100014
$ws.Range("K132").Value = 16880.5005
$ws.Range("L132").Value = 300042
$ws.Range("M132").Value = -14350.5005
$ws.Range("N132").Value = -305102

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H94").Value = 1010.3714
$ws.Range("I94").Value = 686.55554
$ws.Range("J94").Value = 1353.2354
$ws.Range("K94").Value = 686.55554
$ws.Range("L94").Value = 1353.2354
$ws.Range("M94").Value = -235.55554
$ws.Range("N94").Value = -2255.2354

$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H134").Value = 1265.5
$ws.Range("I134").Value = 883
$ws.Range("J134").Value = 3560.5
$ws.Range("K134").Value = 2649
$ws.Range("L134").Value = 10681.5
$ws.Range("M134").Value = -114
$ws.Range("N134").Value = -15751.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1075.7222
$ws.Range("I94").Value = 788.7778
$ws.Range("J94").Value = 1362.6666
$ws.Range("K94").Value = 788.7778
$ws.Range("L94").Value = 1362.6666
$ws.Range("M94").Value = -337.7778
$ws.Range("N94").Value = -2264.6666

$ws.Range("H99").Value = 2307.2222
$ws.Range("I99").Value = 2413.1
$ws.Range("J99").Value = 2174.875
$ws.Range("K99").Value = 2413.1
$ws.Range("L99").Value = 2174.875
$ws.Range("M99").Value = -915.0999999999999
$ws.Range("N99").Value = -5170.875

$ws.Range("H126").Value = 2307.2222
$ws.Range("I126").Value = 2413.1
$ws.Range("J126").Value = 2174.875
$ws.Range("K126").Value = 7239.299999999999
$ws.Range("L126").Value = 6524.625
$ws.Range("M126").Value = -4769.299999999999
$ws.Range("N126").Value = -11464.625

$ws.Range("H132").Value = 66676510
$ws.Range("I132").Value = 125014780
$ws.Range("J132").Value = 4205.5713
$ws.Range("K132").Value = 375044340
$ws.Range("L132").Value = 12616.7139
$ws.Range("M132").Value = -375041810
$ws.Range("N132").Value = -17676.7139

$ws.Range("H134").Value = 1884.9744
$ws.Range("I134").Value = 2061.8462
$ws.Range("J134").Value = 1531.2307
$ws.Range("K134").Value = 6185.5386
$ws.Range("L134").Value = 4593.6921
$ws.Range("M134").Value = -3650.5386
$ws.Range("N134").Value = -9663.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 595.58826
$ws.Range("I2").Value = 918.9
$ws.Range("J2").Value = 133.71428
$ws.Range("K2").Value = 5513.4
$ws.Range("L2").Value = 802.28568
$ws.Range("M2").Value = -5400.4
$ws.Range("N2").Value = -1028.28568

$ws.Range("H17").Value = 5845.7334
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 6191.857
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 18575.571
$ws.Range("M17").Value = -2831
$ws.Range("N17").Value = -18913.571

$ws.Range("H42").Value = 3000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 3000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 9000
$ws.Range("N42").Value = -10068

$ws.Range("H109").Value = 2285.4
$ws.Range("I109").Value = 513.5
$ws.Range("J109").Value = 3466.6667
$ws.Range("K109").Value = 1540.5
$ws.Range("L109").Value = 10400.0001
$ws.Range("M109").Value = -500.5
$ws.Range("N109").Value = -12480.0001

$ws.Range("H120").Value = 12405
$ws.Range("I120").Value = 10518
$ws.Range("J120").Value = 13663
$ws.Range("K120").Value = 31554
$ws.Range("L120").Value = 40989
$ws.Range("M120").Value = -26716
$ws.Range("N120").Value = -50665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws.Range("H122").Value = 2453.9285
$ws.Range("I122").Value = 1949.7273
$ws.Range("J122").Value = 4302.6665
$ws.Range("K122").Value = 5849.1819
$ws.Range("L122").Value = 12907.9995
$ws.Range("M122").Value = -3399.1819
$ws.Range("N122").Value = -17807.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7290.2
$ws.Range("I22").Value = 401
$ws.Range("J22").Value = 8350.076999999999
$ws.Range("K22").Value = 401
$ws.Range("L22").Value = 8350.076999999999
$ws.Range("M22").Value = -106
$ws.Range("N22").Value = -8940.076999999999

$ws.Range("H27").Value = 7290.2
$ws.Range("I27").Value = 401
$ws.Range("J27").Value = 8350.076999999999
$ws.Range("K27").Value = 401
$ws.Range("L27").Value = 8350.076999999999
$ws.Range("M27").Value = -294
$ws.Range("N27").Value = -8564.076999999999

$ws.Range("H62").Value = 43474.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 43474.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 43474.5
$ws.Range("N62").Value = -44722.5

$ws.Range("H64").Value = 31757.143
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 31757.143
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 31757.143
$ws.Range("N64").Value = -32207.143

$ws.Range("H65").Value = 43474.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 43474.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 130423.5
$ws.Range("N65").Value = -136663.5

$ws.Range("H67").Value = 31757.143
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 31757.143
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 31757.143
$ws.Range("N67").Value = -33317.143

$ws.Range("H100").Value = 1413.6342
$ws.Range("I100").Value = 1306.3667
$ws.Range("J100").Value = 1706.1818
$ws.Range("K100").Value = 1306.3667
$ws.Range("L100").Value = 1706.1818
$ws.Range("M100").Value = -765.3667
$ws.Range("N100").Value = -2788.1818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17680
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 17680
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 17680
$ws.Range("N54").Value = -18720

$ws.Range("H62").Value = 2414.4285
$ws.Range("I62").Value = 2157.4285
$ws.Range("J62").Value = 2671.4285
$ws.Range("K62").Value = 2157.4285
$ws.Range("L62").Value = 2671.4285
$ws.Range("M62").Value = -1533.4285
$ws.Range("N62").Value = -3919.4285

$ws.Range("H63").Value = 45000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 45000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 45000
$ws.Range("N63").Value = -46248

$ws.Range("H65").Value = 2414.4285
$ws.Range("I65").Value = 2157.4285
$ws.Range("J65").Value = 2671.4285
$ws.Range("K65").Value = 10787.1425
$ws.Range("L65").Value = 13357.1425
$ws.Range("M65").Value = -7667.1425
$ws.Range("N65").Value = -19597.1425

$ws.Range("H66").Value = 45000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 45000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 135000
$ws.Range("N66").Value = -141240

$ws.Range("H132").Value = 3944.889
$ws.Range("I132").Value = 4772.2905
$ws.Range("J132").Value = 2112.7856
$ws.Range("K132").Value = 14316.8715
$ws.Range("L132").Value = 6338.3568
$ws.Range("M132").Value = -11786.8715
$ws.Range("N132").Value = -11398.3568
